$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Testcodeword5"
$ws.Range("A6").Value = "Testcodeword6"
$ws.Range("A7").Value = "Testcodeword7"
$ws.Range("A8").Value = "Testcodeword8"
$ws.Range("A9").Value = "Testcodeword9"
$ws.Range("A10").Value = "Testcodeword10"
$ws.Range("A11").Value = "Testcodeword11"
$ws.Range("A12").Value = "Testcodeword12"

$ws.Range("G8").Select()
